$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new data rows appended after the existing last row (203),
# for the "Macroferia Regional de Talca" / "Sandia" weekly update.
$rows = @(
    @{ Row = 204; Quality = "Extra";    Volumen = 5000; Precio = 2800 },
    @{ Row = 205; Quality = "Primera";  Volumen = 5000; Precio = 2300 },
    @{ Row = 206; Quality = "Segunda";  Volumen = 5000; Precio = 1800 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"

    # 44911 = 2022-12-16 as an Excel serial date; set the raw serial then the
    # number format so the cell reuses the workbook's existing date style
    # instead of a new auto-generated one.
    $ws.Cells.Item($row, 4).Value = 44911
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = 100112028
    $ws.Cells.Item($row, 7).Value = "Sandia"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.Quality
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Precio
    $ws.Cells.Item($row, 12).Value = $r.Precio
    $ws.Cells.Item($row, 13).Value = $r.Precio
    $ws.Cells.Item($row, 14).Value = "$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
